$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "GRT-USD"
